$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.168866515159607
$ws.Range("B1").Value = 2.441378831863403
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.364012956619263
$ws.Range("E1").Value = 1.236058592796326
